$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# (updates every cell that currently shows the old status so the shared
# string is fully replaced, on every localization-status sheet)

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"
$ws1.Range("E4").Value = "In Translation"
$ws1.Range("F4").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("C4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("C4").Value = "In Translation"

# --- Shrink the Status columns to fit the new, shorter text ---
# Overview sheet: zh-cn (E) and de-de (F) status columns
$ws1.Range("E1").ColumnWidth = 12.5
$ws1.Range("F1").ColumnWidth = 12.5

# zh-cn sheet: Status column (C)
$ws2.Range("C1").ColumnWidth = 12.5

# de-de sheet: Status column (C)
$ws3.Range("C1").ColumnWidth = 12.5
